$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits right
#    after the word "Bahan" (2.2 Alat dan Bahan heading).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Change "6 November 2018" -> "9 November 2018", but keep the
#    resulting text split across two runs ("9" and " November 2018"),
#    matching how Word naturally splits a run when only part of it is
#    retyped.
# ------------------------------------------------------------------
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Text = "6 November 2018"
$found = $find.Find.Execute()

if ($found) {
    $dateStart = $find.Start

    # Temporary bookmarks anchor the two run boundaries we want to
    # keep so that the engine does not coalesce them back together
    # with neighbouring runs that share identical formatting.
    $anchorBefore = $d.Range($dateStart, $dateStart)
    $d.Bookmarks.Add("ZZZ_ANCHOR_BEFORE", $anchorBefore) | Out-Null

    $anchorMid = $d.Range($dateStart + 1, $dateStart + 1)
    $d.Bookmarks.Add("ZZZ_ANCHOR_MID", $anchorMid) | Out-Null

    $digit = $d.Range($dateStart, $dateStart + 1)
    $digit.Text = "9"

    $d.Bookmarks("ZZZ_ANCHOR_BEFORE").Delete()
    $d.Bookmarks("ZZZ_ANCHOR_MID").Delete()
}

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the very end of the document
#    (end of the last paragraph, right before the final paragraph
#    mark), reflecting the last edit position.
#
#    A collapsed range placed exactly on the last character position
#    of a paragraph (immediately before its paragraph mark) is typed
#    in temporarily with a placeholder character so the bookmark can
#    be anchored at a "safe" position, then the placeholder is
#    removed again - the bookmark stays put.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$endPos = $lastPara.Range.End - 1

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("Z")

$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

$placeholderRange = $d.Range($endPos, $endPos + 1)
$placeholderRange.Delete()
